$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.516.58"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "3.229.72"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'604.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'158.11"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.228.20"
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("D11").Value = "'5.67"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -6.56%  "
$ws.Range("D12").Value = "'0.508"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("D14").Value = "'39.03"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "3.761.69"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").Value = "66.620.84"
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("D17").Value = "'7.44"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "3.240.87"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").Value = "'510.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "'15.29"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "'8.08"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "'14.76"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").Value = "'84.80"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").Value = "'9.17"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("E29").Value = "  +5.08%  "
$ws.Range("D30").Value = "'3.00"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.67%  "
$ws.Range("D31").Value = "'7.05"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").Value = "'28.22"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").Value = "'1.19"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("D35").Value = "'6.52"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.0972"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +7.77%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'515.95"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +7.21%  "
$ws.Range("D38").Value = "'56.15"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("E39").Value = "  +18.06%  "
$ws.Range("D40").Value = "'0.0421"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").Value = "'3.04"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.94%  "
$ws.Range("E42").Value = "  +5.88%  "
$ws.Range("D43").Value = "'8.79"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("D46").Value = "2.880.58"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").Value = "'28.59"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("E48").Value = "  +4.60%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "'2.65"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.21%  "
